# Auto-generated edit script: update FFXIV leve-profit market data values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets, per the scheduled
# market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H129").Value = 1526.1
$ws.Range("J129").Value = 1000
$ws.Range("L129").Value = 3000
$ws.Range("N129").Value = -13000
$ws.Range("H132").Value = 45481
$ws.Range("I132").Value = 30084.379
$ws.Range("J132").Value = 116690.375
$ws.Range("K132").Value = 90253.137
$ws.Range("L132").Value = 350071.125
$ws.Range("M132").Value = -87723.137
$ws.Range("N132").Value = -355131.125

$ws = $wb.Worksheets("ARM")
$ws.Range("H58").Value = 37693.332
$ws.Range("J58").Value = 37693.332
$ws.Range("L58").Value = 37693.332
$ws.Range("N58").Value = -38553.332
$ws.Range("H61").Value = 2832.875
$ws.Range("I61").Value = 1962.6
$ws.Range("J61").Value = 3454.5
$ws.Range("K61").Value = 1962.6
$ws.Range("L61").Value = 3454.5
$ws.Range("M61").Value = -1750.6
$ws.Range("N61").Value = -3878.5
$ws.Range("H132").Value = 2530.2917
$ws.Range("I132").Value = 1962.1025
$ws.Range("J132").Value = 4992.4443
$ws.Range("K132").Value = 5886.3075
$ws.Range("L132").Value = 14977.3329
$ws.Range("M132").Value = -3356.3075
$ws.Range("N132").Value = -20037.3329
$ws.Range("H136").Value = 2832.875
$ws.Range("I136").Value = 1962.6
$ws.Range("J136").Value = 3454.5
$ws.Range("K136").Value = 5887.799999999999
$ws.Range("L136").Value = 10363.5
$ws.Range("M136").Value = -3337.799999999999
$ws.Range("N136").Value = -15463.5

$ws = $wb.Worksheets("BSM")
$ws.Range("H20").Value = 6392.522
$ws.Range("J20").Value = 7417.263
$ws.Range("L20").Value = 7417.263
$ws.Range("N20").Value = -7911.263
$ws.Range("H132").Value = 54062.5
$ws.Range("J132").Value = 54062.5
$ws.Range("L132").Value = 54062.5
$ws.Range("N132").Value = -64182.5
$ws.Range("H133").Value = 55498.5
$ws.Range("J133").Value = 55498.5
$ws.Range("L133").Value = 55498.5
$ws.Range("N133").Value = -65618.5

$ws = $wb.Worksheets("CRP")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H132").Value = 51807.83
$ws.Range("I132").Value = 2472.4443
$ws.Range("J132").Value = 132538.45
$ws.Range("K132").Value = 7417.3329
$ws.Range("L132").Value = 397615.35
$ws.Range("M132").Value = -4887.3329
$ws.Range("N132").Value = -402675.35
$ws.Range("H133").Value = 25666.666
$ws.Range("J133").Value = 25666.666
$ws.Range("L133").Value = 25666.666
$ws.Range("N133").Value = -30726.666
$ws.Range("H134").Value = 1278087.6
$ws.Range("I134").Value = 1175018.4
$ws.Range("J134").Value = 1401770.9
$ws.Range("K134").Value = 3525055.2
$ws.Range("L134").Value = 4205312.699999999
$ws.Range("M134").Value = -3522520.2
$ws.Range("N134").Value = -4210382.699999999

$ws = $wb.Worksheets("CUL")
$ws.Range("H18").Value = 25149.834
$ws.Range("I18").Value = 27408.908
$ws.Range("K18").Value = 82226.724
$ws.Range("M18").Value = -82057.724
$ws.Range("H68").Value = 1326.8148
$ws.Range("I68").Value = 1196.25
$ws.Range("K68").Value = 3588.75
$ws.Range("M68").Value = -2777.75
$ws.Range("H71").Value = 1326.8148
$ws.Range("I71").Value = 1196.25
$ws.Range("K71").Value = 10766.25
$ws.Range("M71").Value = -6710.25
$ws.Range("H107").Value = 3462.4775
$ws.Range("J107").Value = 5025.5
$ws.Range("L107").Value = 15076.5
$ws.Range("N107").Value = -18916.5
$ws.Range("H131").Value = 5074.5
$ws.Range("I131").Value = 14697
$ws.Range("J131").Value = 1867
$ws.Range("K131").Value = 44091
$ws.Range("L131").Value = 5601
$ws.Range("M131").Value = -39051
$ws.Range("N131").Value = -15681
$ws.Range("H136").Value = 22729552
$ws.Range("I136").Value = 45456300
$ws.Range("J136").Value = 2803
$ws.Range("K136").Value = 136368900
$ws.Range("L136").Value = 8409
$ws.Range("M136").Value = -136363800
$ws.Range("N136").Value = -18609

$ws = $wb.Worksheets("GSM")
$ws.Range("H70").Value = 5749.9443
$ws.Range("I70").Value = 5700.25
$ws.Range("J70").Value = 6147.5
$ws.Range("K70").Value = 5700.25
$ws.Range("L70").Value = 6147.5
$ws.Range("M70").Value = -5430.25
$ws.Range("N70").Value = -6687.5
$ws.Range("H73").Value = 5749.9443
$ws.Range("I73").Value = 5700.25
$ws.Range("J73").Value = 6147.5
$ws.Range("K73").Value = 5700.25
$ws.Range("L73").Value = 6147.5
$ws.Range("M73").Value = -4764.25
$ws.Range("N73").Value = -8019.5
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H132").Value = 33336286
$ws.Range("I132").Value = 58825732
$ws.Range("K132").Value = 176477196
$ws.Range("M132").Value = -176474666
$ws.Range("H133").Value = 62518
$ws.Range("J133").Value = 62518
$ws.Range("L133").Value = 62518
$ws.Range("N133").Value = -72638
$ws.Range("H135").Value = 70832
$ws.Range("J135").Value = 70832
$ws.Range("L135").Value = 70832
$ws.Range("N135").Value = -80972
$ws.Range("H137").Value = 56533
$ws.Range("J137").Value = 56533
$ws.Range("L137").Value = 56533
$ws.Range("N137").Value = -66733
$ws.Range("H138").Value = 44766.125
$ws.Range("J138").Value = 44766.125
$ws.Range("L138").Value = 44766.125
$ws.Range("N138").Value = -55046.125
$ws.Range("H139").Value = 32255
$ws.Range("J139").Value = 32255
$ws.Range("L139").Value = 32255
$ws.Range("N139").Value = -42535

$ws = $wb.Worksheets("LTW")
$ws.Range("H132").Value = 5597.8887
$ws.Range("I132").Value = 4199
$ws.Range("J132").Value = 5997.5713
$ws.Range("K132").Value = 12597
$ws.Range("L132").Value = 17992.7139
$ws.Range("M132").Value = -10067
$ws.Range("N132").Value = -23052.7139

$ws = $wb.Worksheets("WVR")
$ws.Range("H81").Value = 2510
$ws.Range("I81").Value = 1800
$ws.Range("J81").Value = 3575
$ws.Range("K81").Value = 3600
$ws.Range("L81").Value = 7150
$ws.Range("M81").Value = -2539
$ws.Range("N81").Value = -9272
$ws.Range("H84").Value = 2510
$ws.Range("I84").Value = 1800
$ws.Range("J84").Value = 3575
$ws.Range("K84").Value = 18000
$ws.Range("L84").Value = 35750
$ws.Range("M84").Value = -12696
$ws.Range("N84").Value = -46358
$ws.Range("H132").Value = 3347849.5
$ws.Range("I132").Value = 7248841.5
$ws.Range("J132").Value = 4142.2856
$ws.Range("K132").Value = 21746524.5
$ws.Range("L132").Value = 12426.8568
$ws.Range("M132").Value = -17746524.5
$ws.Range("N132").Value = -17486.8568

